# Insert a new data row at row 342 (pushes existing rows 342:430 down to 343:431)
# and populate it with a new price record (copying category fields that repeat
# throughout the dataset, with a new date and new price/volume figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 342.
$ws.Rows("342:342").Insert()

# Populate the newly inserted row 342 with the new record's data.
$ws.Cells.Item(342, 1).Value = 5
$ws.Cells.Item(342, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(342, 3).Value = "Maule"
$ws.Cells.Item(342, 4).Value = 44722
$ws.Cells.Item(342, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(342, 5).Value = 7
$ws.Cells.Item(342, 6).Value = 100112043
$ws.Cells.Item(342, 7).Value = "Pepino ensalada"
$ws.Cells.Item(342, 8).Value = "Sin especificar"
$ws.Cells.Item(342, 9).Value = "Primera"
$ws.Cells.Item(342, 10).Value = 300
$ws.Cells.Item(342, 11).Value = 20000
$ws.Cells.Item(342, 12).Value = 20000
$ws.Cells.Item(342, 13).Value = 20000
$ws.Cells.Item(342, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(342, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(342, 16).Value = 333
$ws.Cells.Item(342, 17).Value = 60
$ws.Cells.Item(342, 18).Value = "Hortaliza"
